# Update the "Koefisien Variasi" (J) column values for rows 2-28.
# Rows 2-11 keep their original values; rows 12-28 are the ones whose
# values actually change (a cyclic re-shuffle of the J12:J28 block).
# The paste also drops the explicit "0.000" number format those cells
# had, reverting them to the workbook's default (General) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newJValues = @(
    3.079,
    0.6,
    1.596,
    7.23,
    0.643,
    4.183,
    0.352,
    2.114,
    0.538,
    1.066,
    2.671,
    0.203,
    0.255,
    3.161,
    0.701,
    0.548,
    3.954,
    2.092,
    1.002,
    0.739,
    1.149,
    0.05,
    1.398,
    0.355,
    32.552,
    3.279,
    0.481
)

# Use a real .NET 2D object array (rows x cols) rather than a PowerShell
# jagged array literal -- the COM bridge here only marshals a true
# System.Object[,] correctly when assigning to a multi-cell Range.Value.
$rowCount = $newJValues.Length
$data = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i, 0] = $newJValues[$i]
}

$targetRange = $ws.Range("J2:J28")
$targetRange.Value = $data

# The incoming paste carried no explicit number format, so the cells
# fall back to the default "Normal" style (clears the previous "0.000"
# custom format, i.e. drops the s="2" style index).
$targetRange.Style = "Normal"

# Reflect the new selection/scroll position left in the sheet after the edit.
$ws.Range("I4").Select()
